# [external commands] - [`tail(id,file)`]: simulate the *NIX tail command.
#
# This script reproduces, via the Excel object model, the edits captured in the
# OOXML diff:
#   1. The hidden "#system" sheet gains a new "tail(id,file)" entry in the
#      "external" commands column (I), so the "external" named range grows
#      by one row ($I$2:$I$4 -> $I$2:$I$5).
#   2. The hidden "#system" sheet gains a new "assertTextNotContains(locator,text)"
#      entry in the "web" commands column (Y), inserted right after the header
#      (at the top of the existing list, row 39) which pushes every other
#      "web" command down by one row, so the "web" named range grows by one
#      row ($Y$2:$Y$127 -> $Y$2:$Y$128).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1. "external" commands: add tail(id,file) -----------------------------
# Column I currently holds rows 2-4 (header in row 1). Append the new
# command right after the last populated row.
$ws.Range("I5").Value2 = "tail(id,file)"

# --- 2. "web" commands: add assertTextNotContains(locator,text) ------------
# Column Y currently holds rows 2-127 (header in row 1, data starting row 2).
# Insert a brand new row at the very top of the data block (row 39, which is
# where the new entry lands in the authoritative edit) so every following
# "web" command value shifts down by one row.
$ws.Rows.Item(39).Insert()
$ws.Range("Y39").Value2 = "assertTextNotContains(locator,text)"

# --- 3. Keep the named ranges in sync with the new list sizes ---------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "external") {
        $n.RefersTo = "='#system'!`$I`$2:`$I`$5"
    }
    elseif ($n.Name -eq "web") {
        $n.RefersTo = "='#system'!`$Y`$2:`$Y`$128"
    }
}
